$wb = $excel.ActiveWorkbook
$wsLogin = $wb.Worksheets.Item("Login")
$wsGroups = $wb.Worksheets.Item("Groups")

# Update the expected-result text for the "userIncorrectLogin" test case (Login!D4)
$wsLogin.Range("D4").Value = "Sorry, something terrible happened to server."

# Update the expected group name for the "addGroupBtnNameTest" test case (Groups!D4)
$wsGroups.Range("D4").Value = "OstrTestGroup6"

# Move the selection on Login to D8 (this sheet is no longer the active tab)
$wsLogin.Range("D8").Select()

# Finally select Groups!E5, making Groups the active sheet/tab
$wsGroups.Range("E5").Select()
